# Update the title date line
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-01-15 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-16 Thursday", 2)

# Update the multiplication table cells by direct (row, column) addressing
# so that values which become duplicates of other (soon-to-change) cells
# don't get clobbered by a later Find/Replace pass.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="86×88=7568"},
    @{Row=1;  Col=2; Text="91×52=4732"},
    @{Row=1;  Col=3; Text="21×84=1764"},
    @{Row=1;  Col=4; Text="52×26=1352"},
    @{Row=1;  Col=5; Text="17×11=187"},

    @{Row=5;  Col=1; Text="68×86=5848"},
    @{Row=5;  Col=2; Text="57×34=1938"},
    @{Row=5;  Col=3; Text="36×52=1872"},
    @{Row=5;  Col=4; Text="62×77=4774"},
    @{Row=5;  Col=5; Text="64×20=1280"},

    @{Row=10; Col=1; Text="26×15=390"},
    @{Row=10; Col=2; Text="76×52=3952"},
    @{Row=10; Col=3; Text="80×66=5280"},
    @{Row=10; Col=4; Text="87×36=3132"},
    @{Row=10; Col=5; Text="36×76=2736"},

    @{Row=15; Col=1; Text="96×63=6048"},
    @{Row=15; Col=2; Text="88×71=6248"},
    @{Row=15; Col=3; Text="61×18=1098"},
    @{Row=15; Col=4; Text="13×16=208"},
    @{Row=15; Col=5; Text="49×53=2597"},

    @{Row=20; Col=1; Text="50×74=3700"},
    @{Row=20; Col=2; Text="14×83=1162"},
    @{Row=20; Col=3; Text="72×88=6336"},
    @{Row=20; Col=4; Text="63×40=2520"},
    @{Row=20; Col=5; Text="86×92=7912"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
